# Issue #5: "property aircraft done"
#
# The "建物" (building) sheet's property_category column (column I) was
# left set to "land" (copy/paste leftover from the 土地 sheet). Correct
# every data row so it reads "building" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)  # column I = property_category
    if ($cell.Value() -eq "land") {
        $cell.Value = "building"
    }
}
